$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 17 (weekly update) — this shifts the existing rows
# 17-25 down to 18-26, preserving their data and formatting untouched.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with this week's data for
# "Femacal de La Calera" / Perejil.
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "Femacal de La Calera"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44669
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 100112044
$ws.Range("G17").Value = "Perejil"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 92
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 2755
$ws.Range("N17").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O17").Value = "Provincia de Quillota"
$ws.Range("P17").Value = 918
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = "Hortaliza"
